# Removing less than USD 5 price from extrapolation calibration because it is just a noise.
# This recalculates the extrapolation-derived columns (ABSM1_RN, M1_RN, CM2_RN, CMN3_RN, CMN4_RN)
# for the affected rows once the sub-$5 priced noisy input is dropped from the calibration set.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (ExpiryDate 2025-12-26)
$ws.Cells.Item(5, 4).Value  = 117150.2461360985
$ws.Cells.Item(5, 5).Value  = -0.006123160580031033
$ws.Cells.Item(5, 6).Value  = 0.2050107471725738
$ws.Cells.Item(5, 7).Value  = -1.433875363915494
$ws.Cells.Item(5, 8).Value  = 12.54071655823056

# Row 8 (ExpiryDate 2026-05-29)
$ws.Cells.Item(8, 4).Value  = 118916.577262021
$ws.Cells.Item(8, 5).Value  = -0.0251494183698092
$ws.Cells.Item(8, 6).Value  = 0.1873896576901241
$ws.Cells.Item(8, 7).Value  = -1.249267059359946
$ws.Cells.Item(8, 8).Value  = 9.107047764240379

# Row 9 (ExpiryDate 2026-06-26)
$ws.Cells.Item(9, 4).Value  = 120416.6278506177
$ws.Cells.Item(9, 5).Value  = -0.05699414791304799
$ws.Cells.Item(9, 6).Value  = 0.3172890691105682
$ws.Cells.Item(9, 7).Value  = -1.872573934556988
$ws.Cells.Item(9, 8).Value  = 11.80292490804769

# Row 10 (ExpiryDate 2026-09-25)
$ws.Cells.Item(10, 4).Value = 122002.9942314485
$ws.Cells.Item(10, 5).Value = -0.09392150898699426
$ws.Cells.Item(10, 6).Value = 0.4010852999498444
$ws.Cells.Item(10, 7).Value = -1.866035533514716
$ws.Cells.Item(10, 8).Value = 9.817061849021904

# Row 11 (ExpiryDate 2026-12-24)
$ws.Cells.Item(11, 4).Value = 123946.6238077755
$ws.Cells.Item(11, 5).Value = -0.1595794594908676
$ws.Cells.Item(11, 6).Value = 0.6667345552668406
$ws.Cells.Item(11, 7).Value = -2.4897779792289
$ws.Cells.Item(11, 8).Value = 11.88858874292346

# Row 13 (ExpiryDate 2025-09-16)
$ws.Cells.Item(13, 4).Value = 114907.3229185066
$ws.Cells.Item(13, 5).Value = 0.03426577685721082
$ws.Cells.Item(13, 6).Value = 0.1072864368906373
$ws.Cells.Item(13, 7).Value = -0.6763856222028046
$ws.Cells.Item(13, 8).Value = 6.83385722220373

# Row 14 (ExpiryDate 2025-09-17)
$ws.Cells.Item(14, 4).Value = 114912.6081470953
$ws.Cells.Item(14, 5).Value = 0.02997203823671524
$ws.Cells.Item(14, 6).Value = 0.1158857749884572
$ws.Cells.Item(14, 7).Value = -0.6877102909796046
$ws.Cells.Item(14, 8).Value = 6.931938019205041

# Row 18 (ExpiryDate 2025-10-03)
$ws.Cells.Item(18, 4).Value = 115846.3616182638
$ws.Cells.Item(18, 5).Value = 0.03460914777348335
$ws.Cells.Item(18, 6).Value = 0.1410538038394651
$ws.Cells.Item(18, 7).Value = -0.746881217033992
$ws.Cells.Item(18, 8).Value = 6.961375465380724
